$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Coupling Parameters")

# Update End Year value from 2025 to 2030
$ws.Range("B3").Value = 2030

# Update the active selection to C3
$ws.Activate()
$ws.Range("C3").Select()
